$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Converted Data")

# Row 7 (Weights) changes: X7 and AH7 drop from 0.5 to 0, AN7 (weight sum) 13 -> 12
$ws.Range("X7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AN7").Value = 12

# Recomputed LockdownEffectiveness (column AN) values for existing data rows 20-221
$anUpdates = @{
    20 = 0.097222222225
    21 = 0.097222222225
    22 = 0.097222222225
    23 = 0.097222222225
    24 = 0.1805555555583333
    25 = 0.1805555555583333
    26 = 0.1805555555583333
    27 = 0.2638888888916667
    28 = 0.2638888888916667
    29 = 0.2638888888916667
    30 = 0.2638888888916667
    31 = 0.2638888888916667
    32 = 0.8194444444916664
    33 = 0.8194444444916664
    34 = 0.8194444444916664
    35 = 0.8194444444916664
    36 = 0.8194444444916664
    37 = 0.8194444444916664
    38 = 0.8194444444916664
    39 = 0.8194444444916664
    40 = 0.8194444444916664
    41 = 0.8194444444916664
    42 = 0.8194444444916664
    43 = 0.8194444444916664
    44 = 0.8194444444916664
    45 = 0.8194444444916664
    46 = 0.833333333383333
    47 = 0.833333333383333
    48 = 0.833333333383333
    49 = 0.833333333383333
    50 = 0.833333333383333
    51 = 0.833333333383333
    52 = 0.833333333383333
    53 = 0.833333333383333
    54 = 0.833333333383333
    55 = 0.833333333383333
    56 = 0.833333333383333
    57 = 0.833333333383333
    58 = 0.833333333383333
    59 = 0.833333333383333
    60 = 0.833333333383333
    61 = 0.833333333383333
    62 = 0.833333333383333
    63 = 0.833333333383333
    64 = 0.833333333383333
    65 = 0.833333333383333
    66 = 0.833333333383333
    67 = 0.833333333383333
    68 = 0.833333333383333
    69 = 0.833333333383333
    70 = 0.833333333383333
    71 = 0.833333333383333
    72 = 0.833333333383333
    73 = 0.833333333383333
    74 = 0.7500000000499997
    75 = 0.7500000000499997
    76 = 0.7500000000499997
    77 = 0.7500000000499997
    78 = 0.7500000000499997
    79 = 0.7500000000499997
    80 = 0.7500000000499997
    81 = 0.7500000000499997
    82 = 0.7500000000499997
    83 = 0.7500000000499997
    84 = 0.7500000000499997
    85 = 0.7500000000499997
    86 = 0.7500000000499997
    87 = 0.7500000000499997
    88 = 0.6666666667166669
    89 = 0.6666666667166669
    90 = 0.6666666667166669
    91 = 0.6666666667166669
    92 = 0.6666666667166669
    93 = 0.5555555556083336
    94 = 0.5555555556083336
    95 = 0.5555555556083336
    96 = 0.5555555556083336
    97 = 0.5555555556083336
    98 = 0.5555555556083336
    99 = 0.5555555556083336
    100 = 0.5555555556083336
    101 = 0.5555555556083336
    102 = 0.4166666667083334
    103 = 0.4166666667083334
    104 = 0.4166666667083334
    105 = 0.4166666667083334
    106 = 0.4166666667083334
    107 = 0.4166666667083334
    108 = 0.4166666667083334
    109 = 0.4166666667083334
    110 = 0.4166666667083334
    111 = 0.4166666667083334
    112 = 0.4166666667083334
    113 = 0.2361111111249999
    114 = 0.2361111111249999
    115 = 0.2361111111249999
    116 = 0.2361111111249999
    117 = 0.2361111111249999
    118 = 0.2361111111249999
    119 = 0.2361111111249999
    120 = 0.2361111111249999
    121 = 0.2361111111249999
    122 = 0.2361111111249999
    123 = 0.2361111111249999
    124 = 0.2361111111249999
    125 = 0.2361111111249999
    126 = 0.2361111111249999
    127 = 0.2361111111249999
    128 = 0.2361111111249999
    129 = 0.2361111111249999
    130 = 0.2361111111249999
    131 = 0.2361111111249999
    132 = 0.2361111111249999
    133 = 0.2361111111249999
    134 = 0.1805555555583333
    135 = 0.1805555555583333
    136 = 0.1805555555583333
    137 = 0.1805555555583333
    138 = 0.1805555555583333
    139 = 0.1805555555583333
    140 = 0.1805555555583333
    141 = 0.1805555555583333
    142 = 0.1805555555583333
    143 = 0.1805555555583333
    144 = 0.1805555555583333
    145 = 0.1805555555583333
    146 = 0.1805555555583333
    147 = 0.1805555555583333
    148 = 0.1805555555583333
    149 = 0.1805555555583333
    150 = 0.1805555555583333
    151 = 0.1805555555583333
    152 = 0.1805555555583333
    153 = 0.1805555555583333
    154 = 0.1805555555583333
    155 = 0.1805555555583333
    156 = 0.1805555555583333
    157 = 0.2638888888916667
    158 = 0.2638888888916667
    159 = 0.2638888888916667
    160 = 0.2638888888916667
    161 = 0.2638888888916667
    162 = 0.2638888888916667
    163 = 0.2638888888916667
    164 = 0.2638888888916667
    165 = 0.2638888888916667
    166 = 0.2638888888916667
    167 = 0.2638888888916667
    168 = 0.2638888888916667
    169 = 0.2638888888916667
    170 = 0.2638888888916667
    171 = 0.2638888888916667
    172 = 0.2638888888916667
    173 = 0.2638888888916667
    174 = 0.2638888888916667
    175 = 0.2638888888916667
    176 = 0.2638888888916667
    177 = 0.2638888888916667
    178 = 0.2638888888916667
    179 = 0.2638888888916667
    180 = 0.2638888888916667
    181 = 0.2638888888916667
    182 = 0.2638888888916667
    183 = 0.2638888888916667
    184 = 0.2638888888916667
    185 = 0.2638888888916667
    186 = 0.2638888888916667
    187 = 0.2638888888916667
    188 = 0.2638888888916667
    189 = 0.2638888888916667
    190 = 0.2638888888916667
    191 = 0.2638888888916667
    192 = 0.2638888888916667
    193 = 0.2638888888916667
    194 = 0.2638888888916667
    195 = 0.2638888888916667
    196 = 0.2638888888916667
    197 = 0.2638888888916667
    198 = 0.2638888888916667
    199 = 0.2638888888916667
    200 = 0.2638888888916667
    201 = 0.2638888888916667
    202 = 0.2638888888916667
    203 = 0.2638888888916667
    204 = 0.2638888888916667
    205 = 0.2638888888916667
    206 = 0.2638888888916667
    207 = 0.2638888888916667
    208 = 0.2638888888916667
    209 = 0.2638888888916667
    210 = 0.2638888888916667
    211 = 0.2638888888916667
    212 = 0.2638888888916667
    213 = 0.2638888888916667
    214 = 0.2638888888916667
    215 = 0.2638888888916667
    216 = 0.2638888888916667
    217 = 0.2638888888916667
    218 = 0.2638888888916667
    219 = 0.2638888888916667
    220 = 0.2638888888916667
    221 = 0.2638888888916667
}
foreach ($row in $anUpdates.Keys) {
    $ws.Cells.Item($row, 40).Value = $anUpdates[$row]
}

# New data rows 222-233 (9/30/2020 .. 10/11/2020), same weight pattern as row 221
$newRowValues = @(0,0,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,1)
$newRowAN = 0.2638888888916667

foreach ($row in 222..233) {
    # Columns B..AM: fixed weight-indicator pattern (same as row 221)
    for ($c = 0; $c -lt $newRowValues.Length; $c++) {
        $ws.Cells.Item($row, 2 + $c).Value = $newRowValues[$c]
    }

    # Column AN: LockdownEffectiveness value
    $ws.Cells.Item($row, 40).Value = $newRowAN
}

# Column A labels, written as text (matching the bold/bordered style used for the rest of column A)
$ws.Range("A221").Copy()
$ws.Range("A222").PasteSpecial(-4122)
$ws.Range("A222").Formula = "=""9/30/2020"""
$ws.Range("A222").Copy()
$ws.Range("A222").PasteSpecial(-4163)
$ws.Range("A222").Copy()
$ws.Range("A223").PasteSpecial(-4122)
$ws.Range("A223").Formula = "=""10/1/2020"""
$ws.Range("A223").Copy()
$ws.Range("A223").PasteSpecial(-4163)
$ws.Range("A223").Copy()
$ws.Range("A224").PasteSpecial(-4122)
$ws.Range("A224").Formula = "=""10/2/2020"""
$ws.Range("A224").Copy()
$ws.Range("A224").PasteSpecial(-4163)
$ws.Range("A224").Copy()
$ws.Range("A225").PasteSpecial(-4122)
$ws.Range("A225").Formula = "=""10/3/2020"""
$ws.Range("A225").Copy()
$ws.Range("A225").PasteSpecial(-4163)
$ws.Range("A225").Copy()
$ws.Range("A226").PasteSpecial(-4122)
$ws.Range("A226").Formula = "=""10/4/2020"""
$ws.Range("A226").Copy()
$ws.Range("A226").PasteSpecial(-4163)
$ws.Range("A226").Copy()
$ws.Range("A227").PasteSpecial(-4122)
$ws.Range("A227").Formula = "=""10/5/2020"""
$ws.Range("A227").Copy()
$ws.Range("A227").PasteSpecial(-4163)
$ws.Range("A227").Copy()
$ws.Range("A228").PasteSpecial(-4122)
$ws.Range("A228").Formula = "=""10/6/2020"""
$ws.Range("A228").Copy()
$ws.Range("A228").PasteSpecial(-4163)
$ws.Range("A228").Copy()
$ws.Range("A229").PasteSpecial(-4122)
$ws.Range("A229").Formula = "=""10/7/2020"""
$ws.Range("A229").Copy()
$ws.Range("A229").PasteSpecial(-4163)
$ws.Range("A229").Copy()
$ws.Range("A230").PasteSpecial(-4122)
$ws.Range("A230").Formula = "=""10/8/2020"""
$ws.Range("A230").Copy()
$ws.Range("A230").PasteSpecial(-4163)
$ws.Range("A230").Copy()
$ws.Range("A231").PasteSpecial(-4122)
$ws.Range("A231").Formula = "=""10/9/2020"""
$ws.Range("A231").Copy()
$ws.Range("A231").PasteSpecial(-4163)
$ws.Range("A231").Copy()
$ws.Range("A232").PasteSpecial(-4122)
$ws.Range("A232").Formula = "=""10/10/2020"""
$ws.Range("A232").Copy()
$ws.Range("A232").PasteSpecial(-4163)
$ws.Range("A232").Copy()
$ws.Range("A233").PasteSpecial(-4122)
$ws.Range("A233").Formula = "=""10/11/2020"""
$ws.Range("A233").Copy()
$ws.Range("A233").PasteSpecial(-4163)

$excel.CutCopyMode = 0
